$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Impressoa Hp deskjet"
$ws.Range("B5").Value = 1520

$ws.Range("A6").Value = "TV samsung"
$ws.Range("B6").Value = 30500
